# Update TPM-derived NATMI edge metrics (Mfng-Notch1) with refreshed values.
# Columns A-F, K, L are unchanged; columns G-J, M-T carry the new TPM-based
# expression / specificity / edge-weight values for rows 2-26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 26.56908266666666
$ws.Cells.Item(2, 8).Value = 79.70724799999999
$ws.Cells.Item(2, 9).Value = 0.7506383589294218
$ws.Cells.Item(2, 10).Value = 0.7506383589294219
$ws.Cells.Item(2, 13).Value = 68.637375
$ws.Cells.Item(2, 14).Value = 205.912125
$ws.Cells.Item(2, 15).Value = 0.5415701538216162
$ws.Cells.Item(2, 16).Value = 0.5415701538216162
$ws.Cells.Item(2, 17).Value = 1823.632090398
$ws.Cells.Item(2, 18).Value = 16412.688813582
$ws.Cells.Item(2, 19).Value = 0.4065233315098125
$ws.Cells.Item(2, 20).Value = 0.4065233315098126
$ws.Cells.Item(3, 7).Value = 26.56908266666666
$ws.Cells.Item(3, 8).Value = 79.70724799999999
$ws.Cells.Item(3, 9).Value = 0.7506383589294218
$ws.Cells.Item(3, 10).Value = 0.7506383589294219
$ws.Cells.Item(3, 15).Value = 0.08718851262838957
$ws.Cells.Item(3, 16).Value = 0.08718851262838957
$ws.Cells.Item(3, 17).Value = 293.5903472915057
$ws.Cells.Item(3, 18).Value = 2642.313125623552
$ws.Cells.Item(3, 19).Value = 0.06544704203687152
$ws.Cells.Item(3, 20).Value = 0.06544704203687153
$ws.Cells.Item(4, 7).Value = 26.56908266666666
$ws.Cells.Item(4, 8).Value = 79.70724799999999
$ws.Cells.Item(4, 9).Value = 0.7506383589294218
$ws.Cells.Item(4, 10).Value = 0.7506383589294219
$ws.Cells.Item(4, 13).Value = 16.21089566666667
$ws.Cells.Item(4, 14).Value = 48.632687
$ws.Cells.Item(4, 15).Value = 0.1279089892319285
$ws.Cells.Item(4, 16).Value = 0.1279089892319285
$ws.Cells.Item(4, 17).Value = 430.7086270683751
$ws.Cells.Item(4, 18).Value = 3876.377643615376
$ws.Cells.Item(4, 19).Value = 0.0960133937693759
$ws.Cells.Item(4, 20).Value = 0.09601339376937591
$ws.Cells.Item(5, 7).Value = 26.56908266666666
$ws.Cells.Item(5, 8).Value = 79.70724799999999
$ws.Cells.Item(5, 9).Value = 0.7506383589294218
$ws.Cells.Item(5, 10).Value = 0.7506383589294219
$ws.Cells.Item(5, 13).Value = 20.32546233333333
$ws.Cells.Item(5, 14).Value = 60.976387
$ws.Cells.Item(5, 15).Value = 0.1603741949973873
$ws.Cells.Item(5, 16).Value = 0.1603741949973873
$ws.Cells.Item(5, 17).Value = 540.0288889725529
$ws.Cells.Item(5, 18).Value = 4860.260000752975
$ws.Cells.Item(5, 19).Value = 0.1203830225474659
$ws.Cells.Item(5, 20).Value = 0.1203830225474659
$ws.Cells.Item(6, 7).Value = 26.56908266666666
$ws.Cells.Item(6, 8).Value = 79.70724799999999
$ws.Cells.Item(6, 9).Value = 0.7506383589294218
$ws.Cells.Item(6, 10).Value = 0.7506383589294219
$ws.Cells.Item(6, 13).Value = 10.513928
$ws.Cells.Item(6, 14).Value = 31.541784
$ws.Cells.Item(6, 15).Value = 0.08295814932067838
$ws.Cells.Item(6, 16).Value = 0.08295814932067838
$ws.Cells.Item(6, 17).Value = 279.3454221833813
$ws.Cells.Item(6, 18).Value = 2514.108799650432
$ws.Cells.Item(6, 19).Value = 0.06227156906589595
$ws.Cells.Item(6, 20).Value = 0.06227156906589595
$ws.Cells.Item(7, 9).Value = 0.004944072121179194
$ws.Cells.Item(7, 10).Value = 0.004944072121179195
$ws.Cells.Item(7, 13).Value = 68.637375
$ws.Cells.Item(7, 14).Value = 205.912125
$ws.Cells.Item(7, 15).Value = 0.5415701538216162
$ws.Cells.Item(7, 16).Value = 0.5415701538216162
$ws.Cells.Item(7, 17).Value = 12.011334712875
$ws.Cells.Item(7, 18).Value = 108.102012415875
$ws.Cells.Item(7, 19).Value = 0.00267756189917218
$ws.Cells.Item(7, 20).Value = 0.002677561899172181
$ws.Cells.Item(8, 9).Value = 0.004944072121179194
$ws.Cells.Item(8, 10).Value = 0.004944072121179195
$ws.Cells.Item(8, 15).Value = 0.08718851262838957
$ws.Cells.Item(8, 16).Value = 0.08718851262838957
$ws.Cells.Item(8, 19).Value = 0.000431066294573101
$ws.Cells.Item(8, 20).Value = 0.0004310662945731011
$ws.Cells.Item(9, 9).Value = 0.004944072121179194
$ws.Cells.Item(9, 10).Value = 0.004944072121179195
$ws.Cells.Item(9, 13).Value = 16.21089566666667
$ws.Cells.Item(9, 14).Value = 48.632687
$ws.Cells.Item(9, 15).Value = 0.1279089892319285
$ws.Cells.Item(9, 16).Value = 0.1279089892319285
$ws.Cells.Item(9, 17).Value = 2.836858108979667
$ws.Cells.Item(9, 18).Value = 25.531722980817
$ws.Cells.Item(9, 19).Value = 0.0006323912677097874
$ws.Cells.Item(9, 20).Value = 0.0006323912677097875
$ws.Cells.Item(10, 9).Value = 0.004944072121179194
$ws.Cells.Item(10, 10).Value = 0.004944072121179195
$ws.Cells.Item(10, 13).Value = 20.32546233333333
$ws.Cells.Item(10, 14).Value = 60.976387
$ws.Cells.Item(10, 15).Value = 0.1603741949973873
$ws.Cells.Item(10, 16).Value = 0.1603741949973873
$ws.Cells.Item(10, 17).Value = 3.556894931946333
$ws.Cells.Item(10, 18).Value = 32.012054387517
$ws.Cells.Item(10, 19).Value = 0.0007929015864431386
$ws.Cells.Item(10, 20).Value = 0.0007929015864431387
$ws.Cells.Item(11, 9).Value = 0.004944072121179194
$ws.Cells.Item(11, 10).Value = 0.004944072121179195
$ws.Cells.Item(11, 13).Value = 10.513928
$ws.Cells.Item(11, 14).Value = 31.541784
$ws.Cells.Item(11, 15).Value = 0.08295814932067838
$ws.Cells.Item(11, 16).Value = 0.08295814932067838
$ws.Cells.Item(11, 17).Value = 1.839905858216
$ws.Cells.Item(11, 18).Value = 16.559152723944
$ws.Cells.Item(11, 19).Value = 0.0004101510732809867
$ws.Cells.Item(11, 20).Value = 0.0004101510732809868
$ws.Cells.Item(12, 7).Value = 5.094400666666666
$ws.Cells.Item(12, 8).Value = 15.283202
$ws.Cells.Item(12, 9).Value = 0.1439286634067062
$ws.Cells.Item(12, 10).Value = 0.1439286634067062
$ws.Cells.Item(12, 13).Value = 68.637375
$ws.Cells.Item(12, 14).Value = 205.912125
$ws.Cells.Item(12, 15).Value = 0.5415701538216162
$ws.Cells.Item(12, 16).Value = 0.5415701538216162
$ws.Cells.Item(12, 17).Value = 349.66628895825
$ws.Cells.Item(12, 18).Value = 3146.99660062425
$ws.Cells.Item(12, 19).Value = 0.0779474683805095
$ws.Cells.Item(12, 20).Value = 0.07794746838050952
$ws.Cells.Item(13, 7).Value = 5.094400666666666
$ws.Cells.Item(13, 8).Value = 15.283202
$ws.Cells.Item(13, 9).Value = 0.1439286634067062
$ws.Cells.Item(13, 10).Value = 0.1439286634067062
$ws.Cells.Item(13, 15).Value = 0.08718851262838957
$ws.Cells.Item(13, 16).Value = 0.08718851262838957
$ws.Cells.Item(13, 17).Value = 56.29350774858311
$ws.Cells.Item(13, 18).Value = 506.641569737248
$ws.Cells.Item(13, 19).Value = 0.01254892608702284
$ws.Cells.Item(13, 20).Value = 0.01254892608702284
$ws.Cells.Item(14, 7).Value = 5.094400666666666
$ws.Cells.Item(14, 8).Value = 15.283202
$ws.Cells.Item(14, 9).Value = 0.1439286634067062
$ws.Cells.Item(14, 10).Value = 0.1439286634067062
$ws.Cells.Item(14, 13).Value = 16.21089566666667
$ws.Cells.Item(14, 14).Value = 48.632687
$ws.Cells.Item(14, 15).Value = 0.1279089892319285
$ws.Cells.Item(14, 16).Value = 0.1279089892319285
$ws.Cells.Item(14, 17).Value = 82.58479769153045
$ws.Cells.Item(14, 18).Value = 743.263179223774
$ws.Cells.Item(14, 19).Value = 0.01840976985785425
$ws.Cells.Item(14, 20).Value = 0.01840976985785425
$ws.Cells.Item(15, 7).Value = 5.094400666666666
$ws.Cells.Item(15, 8).Value = 15.283202
$ws.Cells.Item(15, 9).Value = 0.1439286634067062
$ws.Cells.Item(15, 10).Value = 0.1439286634067062
$ws.Cells.Item(15, 13).Value = 20.32546233333333
$ws.Cells.Item(15, 14).Value = 60.976387
$ws.Cells.Item(15, 15).Value = 0.1603741949973873
$ws.Cells.Item(15, 16).Value = 0.1603741949973873
$ws.Cells.Item(15, 17).Value = 103.5460488612415
$ws.Cells.Item(15, 18).Value = 931.914439751174
$ws.Cells.Item(15, 19).Value = 0.02308244353090043
$ws.Cells.Item(15, 20).Value = 0.02308244353090043
$ws.Cells.Item(16, 7).Value = 5.094400666666666
$ws.Cells.Item(16, 8).Value = 15.283202
$ws.Cells.Item(16, 9).Value = 0.1439286634067062
$ws.Cells.Item(16, 10).Value = 0.1439286634067062
$ws.Cells.Item(16, 13).Value = 10.513928
$ws.Cells.Item(16, 14).Value = 31.541784
$ws.Cells.Item(16, 15).Value = 0.08295814932067838
$ws.Cells.Item(16, 16).Value = 0.08295814932067838
$ws.Cells.Item(16, 17).Value = 53.56216181248533
$ws.Cells.Item(16, 18).Value = 482.059456312368
$ws.Cells.Item(16, 19).Value = 0.01194005555041919
$ws.Cells.Item(16, 20).Value = 0.01194005555041919
$ws.Cells.Item(17, 7).Value = 0.05498833333333333
$ws.Cells.Item(17, 8).Value = 0.164965
$ws.Cells.Item(17, 9).Value = 0.001553548265532792
$ws.Cells.Item(17, 10).Value = 0.001553548265532792
$ws.Cells.Item(17, 13).Value = 68.637375
$ws.Cells.Item(17, 14).Value = 205.912125
$ws.Cells.Item(17, 15).Value = 0.5415701538216162
$ws.Cells.Item(17, 16).Value = 0.5415701538216162
$ws.Cells.Item(17, 17).Value = 3.774254855625
$ws.Cells.Item(17, 18).Value = 33.968293700625
$ws.Cells.Item(17, 19).Value = 0.000841355373133899
$ws.Cells.Item(17, 20).Value = 0.0008413553731338991
$ws.Cells.Item(18, 7).Value = 0.05498833333333333
$ws.Cells.Item(18, 8).Value = 0.164965
$ws.Cells.Item(18, 9).Value = 0.001553548265532792
$ws.Cells.Item(18, 10).Value = 0.001553548265532792
$ws.Cells.Item(18, 15).Value = 0.08718851262838957
$ws.Cells.Item(18, 16).Value = 0.08718851262838957
$ws.Cells.Item(18, 17).Value = 0.607625189128889
$ws.Cells.Item(18, 18).Value = 5.46862670216
$ws.Cells.Item(18, 19).Value = 0.0001354515625682185
$ws.Cells.Item(18, 20).Value = 0.0001354515625682185
$ws.Cells.Item(19, 7).Value = 0.05498833333333333
$ws.Cells.Item(19, 8).Value = 0.164965
$ws.Cells.Item(19, 9).Value = 0.001553548265532792
$ws.Cells.Item(19, 10).Value = 0.001553548265532792
$ws.Cells.Item(19, 13).Value = 16.21089566666667
$ws.Cells.Item(19, 14).Value = 48.632687
$ws.Cells.Item(19, 15).Value = 0.1279089892319285
$ws.Cells.Item(19, 16).Value = 0.1279089892319285
$ws.Cells.Item(19, 17).Value = 0.8914101345505557
$ws.Cells.Item(19, 18).Value = 8.022691210955001
$ws.Cells.Item(19, 19).Value = 0.000198712788367315
$ws.Cells.Item(19, 20).Value = 0.0001987127883673151
$ws.Cells.Item(20, 7).Value = 0.05498833333333333
$ws.Cells.Item(20, 8).Value = 0.164965
$ws.Cells.Item(20, 9).Value = 0.001553548265532792
$ws.Cells.Item(20, 10).Value = 0.001553548265532792
$ws.Cells.Item(20, 13).Value = 20.32546233333333
$ws.Cells.Item(20, 14).Value = 60.976387
$ws.Cells.Item(20, 15).Value = 0.1603741949973873
$ws.Cells.Item(20, 16).Value = 0.1603741949973873
$ws.Cells.Item(20, 17).Value = 1.117663297939445
$ws.Cells.Item(20, 18).Value = 10.058969681455
$ws.Cells.Item(20, 19).Value = 0.0002491490524744088
$ws.Cells.Item(20, 20).Value = 0.0002491490524744088
$ws.Cells.Item(21, 7).Value = 0.05498833333333333
$ws.Cells.Item(21, 8).Value = 0.164965
$ws.Cells.Item(21, 9).Value = 0.001553548265532792
$ws.Cells.Item(21, 10).Value = 0.001553548265532792
$ws.Cells.Item(21, 13).Value = 10.513928
$ws.Cells.Item(21, 14).Value = 31.541784
$ws.Cells.Item(21, 15).Value = 0.08295814932067838
$ws.Cells.Item(21, 16).Value = 0.08295814932067838
$ws.Cells.Item(21, 17).Value = 0.5781433775066667
$ws.Cells.Item(21, 18).Value = 5.20329039756
$ws.Cells.Item(21, 19).Value = 0.0001288794889889502
$ws.Cells.Item(21, 20).Value = 0.0001288794889889503
$ws.Cells.Item(22, 7).Value = 3.501848333333333
$ws.Cells.Item(22, 8).Value = 10.505545
$ws.Cells.Item(22, 9).Value = 0.09893535727715995
$ws.Cells.Item(22, 10).Value = 0.09893535727715996
$ws.Cells.Item(22, 13).Value = 68.637375
$ws.Cells.Item(22, 14).Value = 205.912125
$ws.Cells.Item(22, 15).Value = 0.5415701538216162
$ws.Cells.Item(22, 16).Value = 0.5415701538216162
$ws.Cells.Item(22, 17).Value = 240.357677248125
$ws.Cells.Item(22, 18).Value = 2163.219095233125
$ws.Cells.Item(22, 19).Value = 0.05358043665898807
$ws.Cells.Item(22, 20).Value = 0.05358043665898808
$ws.Cells.Item(23, 7).Value = 3.501848333333333
$ws.Cells.Item(23, 8).Value = 10.505545
$ws.Cells.Item(23, 9).Value = 0.09893535727715995
$ws.Cells.Item(23, 10).Value = 0.09893535727715996
$ws.Cells.Item(23, 15).Value = 0.08718851262838957
$ws.Cells.Item(23, 16).Value = 0.08718851262838957
$ws.Cells.Item(23, 17).Value = 38.69568555467556
$ws.Cells.Item(23, 18).Value = 348.26116999208
$ws.Cells.Item(23, 19).Value = 0.008626026647353895
$ws.Cells.Item(23, 20).Value = 0.008626026647353896
$ws.Cells.Item(24, 7).Value = 3.501848333333333
$ws.Cells.Item(24, 8).Value = 10.505545
$ws.Cells.Item(24, 9).Value = 0.09893535727715995
$ws.Cells.Item(24, 10).Value = 0.09893535727715996
$ws.Cells.Item(24, 13).Value = 16.21089566666667
$ws.Cells.Item(24, 14).Value = 48.632687
$ws.Cells.Item(24, 15).Value = 0.1279089892319285
$ws.Cells.Item(24, 16).Value = 0.1279089892319285
$ws.Cells.Item(24, 17).Value = 56.76809797215723
$ws.Cells.Item(24, 18).Value = 510.912881749415
$ws.Cells.Item(24, 19).Value = 0.01265472154862125
$ws.Cells.Item(24, 20).Value = 0.01265472154862125
$ws.Cells.Item(25, 7).Value = 3.501848333333333
$ws.Cells.Item(25, 8).Value = 10.505545
$ws.Cells.Item(25, 9).Value = 0.09893535727715995
$ws.Cells.Item(25, 10).Value = 0.09893535727715996
$ws.Cells.Item(25, 13).Value = 20.32546233333333
$ws.Cells.Item(25, 14).Value = 60.976387
$ws.Cells.Item(25, 15).Value = 0.1603741949973873
$ws.Cells.Item(25, 16).Value = 0.1603741949973873
$ws.Cells.Item(25, 17).Value = 71.17668639621277
$ws.Cells.Item(25, 18).Value = 640.590177565915
$ws.Cells.Item(25, 19).Value = 0.01586667828010344
$ws.Cells.Item(25, 20).Value = 0.01586667828010344
$ws.Cells.Item(26, 7).Value = 3.501848333333333
$ws.Cells.Item(26, 8).Value = 10.505545
$ws.Cells.Item(26, 9).Value = 0.09893535727715995
$ws.Cells.Item(26, 10).Value = 0.09893535727715996
$ws.Cells.Item(26, 13).Value = 10.513928
$ws.Cells.Item(26, 14).Value = 31.541784
$ws.Cells.Item(26, 15).Value = 0.08295814932067838
$ws.Cells.Item(26, 16).Value = 0.08295814932067838
$ws.Cells.Item(26, 17).Value = 36.81818124358666
$ws.Cells.Item(26, 18).Value = 279.3454221833813
$ws.Cells.Item(26, 19).Value = 0.0082074941420933
$ws.Cells.Item(26, 20).Value = 0.008207494142093302
$wb.Save()
